$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear existing hyperlinks and content first
$ws.Hyperlinks.Delete()
$ws.Cells.Clear()

# Header row
$ws.Cells.Item(1,1).Value = '取得日時'
$ws.Cells.Item(1,2).Value = 'タイトル'
$ws.Cells.Item(1,3).Value = 'カテゴリ'
$ws.Cells.Item(1,4).Value = '価格'
$ws.Cells.Item(1,5).Value = '締切'
$ws.Cells.Item(1,6).Value = 'URL'
$ws.Cells.Item(1,7).Value = '優先度スコア'
$ws.Cells.Item(1,8).Value = 'スキル概要'

# Row 2
$ws.Cells.Item(2,1).Value = '2026-01-08 18:26:03'
$ws.Cells.Item(2,2).Value = 'EC×AIプロダクト/業務改善リード'
$ws.Cells.Item(2,3).Value = 'システム開発'
$ws.Cells.Item(2,4).Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Cells.Item(2,5).Value = '期限情報なし'
$ws.Cells.Item(2,6).Value = 'https://www.lancers.jp/work/detail/5467702'
$ws.Cells.Item(2,7).Value = 338
$ws.Cells.Item(2,8).Value = '🔥AI,Ai ◇業務改善'
$null = $ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5467702')

# Row 3
$ws.Cells.Item(3,1).Value = '2026-01-08 18:26:03'
$ws.Cells.Item(3,2).Value = '初回 急募 自動カートインツール 購入bot開発のプロフェッショナルを探しています'
$ws.Cells.Item(3,3).Value = 'システム開発'
$ws.Cells.Item(3,4).Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Cells.Item(3,5).Value = '期限情報なし'
$ws.Cells.Item(3,6).Value = 'https://www.lancers.jp/work/detail/5467745'
$ws.Cells.Item(3,7).Value = 235
$ws.Cells.Item(3,8).Value = '★bot ◆ツール,開発'
$null = $ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5467745')

# Row 4
$ws.Cells.Item(4,1).Value = '2026-01-08 18:26:03'
$ws.Cells.Item(4,2).Value = '【急募】Ecommerce開発のシニアデベロッパーを探しています'
$ws.Cells.Item(4,3).Value = 'システム開発'
$ws.Cells.Item(4,4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(4,5).Value = '期限情報なし'
$ws.Cells.Item(4,6).Value = 'https://www.lancers.jp/work/detail/5467698'
$ws.Cells.Item(4,7).Value = 75
$ws.Cells.Item(4,8).Value = '◆開発'
$null = $ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5467698')

# Row 5
$ws.Cells.Item(5,1).Value = '2026-01-08 18:26:03'
$ws.Cells.Item(5,2).Value = '【法人歓迎】プローバステージ制御ソフト開発の見積依頼'
$ws.Cells.Item(5,3).Value = 'システム開発'
$ws.Cells.Item(5,4).Value = '1,000,000 円 ~ 3,000,000 円 / 固定'
$ws.Cells.Item(5,5).Value = '期限情報なし'
$ws.Cells.Item(5,6).Value = 'https://www.lancers.jp/work/detail/5467295'
$ws.Cells.Item(5,7).Value = 75
$ws.Cells.Item(5,8).Value = '◆開発'
$null = $ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5467295')

# Row 6
$ws.Cells.Item(6,1).Value = '2026-01-08 18:26:03'
$ws.Cells.Item(6,2).Value = 'イベントサイトのWeb制作(決済機能付き)依頼'
$ws.Cells.Item(6,3).Value = 'システム開発'
$ws.Cells.Item(6,4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(6,5).Value = '期限情報なし'
$ws.Cells.Item(6,6).Value = 'https://www.lancers.jp/work/detail/5467460'
$ws.Cells.Item(6,7).Value = 38
$ws.Cells.Item(6,8).Value = '◇サイト'
$null = $ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5467460')

# Row 7
$ws.Cells.Item(7,1).Value = '2026-01-08 18:26:03'
$ws.Cells.Item(7,2).Value = '【急募】社内Webアプリの修正・再構築依頼'
$ws.Cells.Item(7,3).Value = 'システム開発'
$ws.Cells.Item(7,4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(7,5).Value = '期限情報なし'
$ws.Cells.Item(7,6).Value = 'https://www.lancers.jp/work/detail/5467384'
$ws.Cells.Item(7,7).Value = 33
$ws.Cells.Item(7,8).Value = '◇アプリ'
$null = $ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5467384')

# Row 8
$ws.Cells.Item(8,1).Value = '2026-01-08 18:26:03'
$ws.Cells.Item(8,2).Value = 'iPhoneのブラウザ要素の書き換えアプリ作成'
$ws.Cells.Item(8,3).Value = 'システム開発'
$ws.Cells.Item(8,4).Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Cells.Item(8,5).Value = '期限情報なし'
$ws.Cells.Item(8,6).Value = 'https://www.lancers.jp/work/detail/5467578'
$ws.Cells.Item(8,7).Value = 30
$ws.Cells.Item(8,8).Value = '◇アプリ'
$null = $ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5467578')

# Row 9
$ws.Cells.Item(9,1).Value = '2026-01-08 18:26:03'
$ws.Cells.Item(9,2).Value = '進行管理およびチームディレクションを担当'
$ws.Cells.Item(9,3).Value = 'システム開発'
$ws.Cells.Item(9,4).Value = '~ 5,000 円 / 固定'
$ws.Cells.Item(9,5).Value = '期限情報なし'
$ws.Cells.Item(9,6).Value = 'https://www.lancers.jp/work/detail/5418064'
$ws.Cells.Item(9,7).Value = 30
$ws.Cells.Item(9,8).Value = '◇管理'
$null = $ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5418064')

# Row 10
$ws.Cells.Item(10,1).Value = '2026-01-08 18:26:03'
$ws.Cells.Item(10,2).Value = '【急募】cloudflare導入の経験者を探しています!'
$ws.Cells.Item(10,3).Value = 'システム開発'
$ws.Cells.Item(10,4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(10,5).Value = '期限情報なし'
$ws.Cells.Item(10,6).Value = 'https://www.lancers.jp/work/detail/5467334'
$ws.Cells.Item(10,7).Value = 18
$null = $ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5467334')

# Row 11
$ws.Cells.Item(11,1).Value = '2026-01-08 18:26:03'
$ws.Cells.Item(11,2).Value = '電気点火装置の回路図作成依頼'
$ws.Cells.Item(11,3).Value = 'システム開発'
$ws.Cells.Item(11,4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(11,5).Value = '期限情報なし'
$ws.Cells.Item(11,6).Value = 'https://www.lancers.jp/work/detail/5466994'
$ws.Cells.Item(11,7).Value = 13
$null = $ws.Hyperlinks.Add($ws.Range("F11"), 'https://www.lancers.jp/work/detail/5466994')

# Row 12
$ws.Cells.Item(12,1).Value = '2026-01-08 18:26:03'
$ws.Cells.Item(12,2).Value = 'リダイレクトでエラーが出てるので修正依頼'
$ws.Cells.Item(12,3).Value = 'システム開発'
$ws.Cells.Item(12,4).Value = '1,000 ~ 5,000 円 / 固定'
$ws.Cells.Item(12,5).Value = '期限情報なし'
$ws.Cells.Item(12,6).Value = 'https://www.lancers.jp/work/detail/5467706'
$ws.Cells.Item(12,7).Value = 10
$null = $ws.Hyperlinks.Add($ws.Range("F12"), 'https://www.lancers.jp/work/detail/5467706')

# Row 13
$ws.Cells.Item(13,1).Value = '2026-01-08 18:26:03'
$ws.Cells.Item(13,2).Value = 'ドメインの移行をして欲しい'
$ws.Cells.Item(13,3).Value = 'システム開発'
$ws.Cells.Item(13,4).Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Cells.Item(13,5).Value = '期限情報なし'
$ws.Cells.Item(13,6).Value = 'https://www.lancers.jp/work/detail/5467598'
$ws.Cells.Item(13,7).Value = 10
$null = $ws.Hyperlinks.Add($ws.Range("F13"), 'https://www.lancers.jp/work/detail/5467598')

# Column widths
$ws.Columns.Item(2).ColumnWidth = 42.1666666666667
$ws.Columns.Item(8).ColumnWidth = 13.1666666666667

$wb.Save()
